$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the project location path (B3)
$ws.Range("B3").Value = "W:\Projects\בהת\154_בית_שמש_מתאר_דרום\קבצי עבודה\תחזיות_דמוגרפיות"

# Update forecast scenario name (B4)
$ws.Range("B4").Value = "full_realization_BS"

# Update v_date (B5)
$ws.Range("B5").Value = 241028

# Update version number factor (B7)
$ws.Range("B7").Value = 1.2

# Update the active selection to B4
$ws.Range("B4").Select()
